$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 22:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8563325
$ws.Range("C4").Value = 42375
$ws.Range("D4").Value = 5564905
$ws.Range("E4").Value = 2771476
$ws.Range("G4").Value = 760
$ws.Range("H4").Value = 226944

# Row 15 - Sudafrica
$ws.Range("B15").Value = 708359
$ws.Range("C15").Value = 2055
$ws.Range("D15").Value = 641706
$ws.Range("E15").Value = 47912
$ws.Range("G15").Value = 85
$ws.Range("H15").Value = 18741

# Row 21 - Alemania
$ws.Range("B21").Value = 391290
$ws.Range("C21").Value = 10392
$ws.Range("E21").Value = 79191
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = 9999

# Row 28 - Israel
$ws.Range("B28").Value = 307335
$ws.Range("C28").Value = 1173
$ws.Range("D28").Value = 285475
$ws.Range("E28").Value = 19569

# Row 31 - Canada
$ws.Range("B31").Value = 205749
$ws.Range("C31").Value = 2061
$ws.Range("D31").Value = 173392
$ws.Range("E31").Value = 22533
$ws.Range("G31").Value = 30
$ws.Range("H31").Value = 9824

# Row 50 - Costa Rica
$ws.Range("B50").Value = 99425
$ws.Range("C50").Value = 1503
$ws.Range("D50").Value = 60738
$ws.Range("E50").Value = 37451
$ws.Range("G50").Value = 14
$ws.Range("H50").Value = 1236

# Row 52 - Suiza
$ws.Range("E52").Value = 35124
$ws.Range("G52").Value = 17
$ws.Range("H52").Value = 2039

# Row 94 - Costa de Marfil
$ws.Range("B94").Value = 20363
$ws.Range("C94").Value = 21
$ws.Range("D94").Value = 20070
$ws.Range("E94").Value = 172

# Row 112 - Guayana Francesa
$ws.Range("B112").Value = 10295
$ws.Range("C112").Value = 27
$ws.Range("E112").Value = 231

# Row 117 - Zimbabue
$ws.Range("B117").Value = 8215
$ws.Range("C117").Value = 28
$ws.Range("D117").Value = 7725
$ws.Range("E117").Value = 254
$ws.Range("G117").Value = 3
$ws.Range("H117").Value = 236

# Row 130 - Trinidad yTobago
$ws.Range("B130").Value = 5392
$ws.Range("C130").Value = 59
$ws.Range("D130").Value = 3822
$ws.Range("E130").Value = 1469
$ws.Range("G130").Value = 3
$ws.Range("H130").Value = 101

# Row 132 - Siria
$ws.Range("B132").Value = 5224
$ws.Range("C132").Value = 44
$ws.Range("D132").Value = 1629
$ws.Range("E132").Value = 3338
$ws.Range("G132").Value = 3
$ws.Range("H132").Value = 257

# Row 141 - Aruba
$ws.Range("B141").Value = 4369
$ws.Range("C141").Value = 14
$ws.Range("D141").Value = 4084
$ws.Range("E141").Value = 250
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 35
